$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Add New Bank Accounts", "FAILED", "chrome"),
    @("Edit The Bank Accounts", "FAILED", "chrome"),
    @("Add New Bank Accounts", "FAILED", "chrome"),
    @("Add New Bank Accounts", "PASSED", "chrome"),
    @("Edit The Bank Accounts", "PASSED", "chrome"),
    @("Delete The Bank Accounts", "PASSED", "chrome")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
